# practica_1.xlsx — "clase 2 de aplicadas 2"
# Replace the old "alumno 1 / alumno 2 / notas / total" mini-table with a
# roster of 3 students + their ages, a header row, a "Suma" total row driven
# by a new named range "edades", and a conditional format that highlights
# ages below 25. The old "rango_notas" named range is left dangling (#REF!)
# since the cells it pointed to are gone.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- Defined names -----------------------------------------------------
# Add the new named range used by the SUM formula below.
[void]$wb.Names.Add("edades", "=Hoja1!`$D`$6:`$D`$8")
# The old named range's target cells are being removed; it now resolves
# to a broken reference, same as the other two legacy names.
$wb.Names.Item("rango_notas").RefersTo = "=Hoja1!#REF!"

# --- Clear the old C3:D6 mini-table -------------------------------------
$ws.Range("C3:D6").Clear()

# --- Write the new roster (order matters for shared-string ids) --------
$ws.Range("C6").Value = "CRIOLLO GARCIA JESSICA VIVIANA"
$ws.Range("D6").Value = 22
$ws.Range("C7").Value = "ECHEVERRIA MOREIRA GIOVANNI JAVIER"
$ws.Range("D7").Value = 27
$ws.Range("C8").Value = "ERAZO LAVAYEN ERIKA LISBETH"
$ws.Range("D8").Value = 23
$ws.Range("C5").Value = "Nombre"
$ws.Range("D5").Value = "Edad"
$ws.Range("C9").Value = "Suma"
$ws.Range("D9").Formula = "=SUM(edades)"

# --- Column C is now much wider (long full names) -----------------------
$ws.Columns.Item(3).ColumnWidth = 37.43

# --- Conditional formatting: green fill when age < 25 ------------------
# Mirrors an edit session where a few candidate rules were tried and
# removed before settling on the final one (matches the 3 leftover dxfs
# in styles.xml, with the surviving rule landing on dxfId index 1).
$rng = $ws.Range("C6:D8")

$fcA = $rng.FormatConditions.Add(2, 0, '$D6<10')
$fcA.Interior.Color = 5287936

$fcB = $rng.FormatConditions.Add(2, 0, '$D6<25')
$fcB.Interior.Color = 5287936

$fcC = $rng.FormatConditions.Add(2, 0, '$D6<30')
$fcC.Interior.Color = 5287936

$rng.FormatConditions.Item(1).Delete()
$rng.FormatConditions.Item(2).Delete()
$rng.FormatConditions.Item(1).SetFirstPriority()

# --- Selection moved to F5 ------------------------------------------------
[void]$ws.Range("F5").Select()
